$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds crypto prices stored as text (values such as "27.960.97" or
# "1.003" must remain literal text, not be reinterpreted as numbers), so we
# prefix with an apostrophe the way Excel expects for forced-text entry.

$ws.Range("D2").Value = "'27.960.97"
$ws.Range("E2").Value = "  +7.12%  "
$ws.Range("D3").Value = "'1.742.76"
$ws.Range("E3").Value = "  +5.45%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'228.06"
$ws.Range("E5").Value = "  +4.21%  "
$ws.Range("D6").Value = "'0.5433"
$ws.Range("E6").Value = "  +3.41%  "
$ws.Range("D8").Value = "'0.2770"
$ws.Range("E8").Value = "  +4.04%  "
$ws.Range("D9").Value = "'0.06746"
$ws.Range("E9").Value = "  +6.24%  "
$ws.Range("E10").Value = "  +5.25%  "
$ws.Range("D11").Value = "'0.07796"
$ws.Range("E11").Value = "  +1.18%  "
$ws.Range("D12").Value = "'4.708"
$ws.Range("E12").Value = "  +2.21%  "
$ws.Range("D13").Value = "'1.732.77"
$ws.Range("E13").Value = "  +2.32%  "
$ws.Range("D14").Value = "'1.982.80"
$ws.Range("E14").Value = "  +5.48%  "
$ws.Range("D15").Value = "'0.5975"
$ws.Range("E15").Value = "  +6.63%  "
$ws.Range("D16").Value = "'0.0₅8395"
$ws.Range("E16").Value = "  +2.49%  "
$ws.Range("D17").Value = "'68.86"
$ws.Range("E17").Value = "  +5.49%  "
$ws.Range("D18").Value = "'27.960.07"
$ws.Range("E18").Value = "  +7.07%  "
$ws.Range("D19").Value = "'224.57"
$ws.Range("E19").Value = "  +17.59%  "
$ws.Range("D20").Value = "'4.844"
$ws.Range("E20").Value = "  +3.11%  "
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("D22").Value = "'10.92"
$ws.Range("E22").Value = "  +5.25%  "
$ws.Range("D23").Value = "'6.243"
$ws.Range("E23").Value = "  +4.30%  "
$ws.Range("D24").Value = "'1.004"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").Value = "'146.34"
$ws.Range("E25").Value = "  +0.32%  "
$ws.Range("D26").Value = "'0.1243"
$ws.Range("E26").Value = "  +3.55%  "
$ws.Range("D27").Value = "'17.28"
$ws.Range("E27").Value = "  +8.55%  "
$ws.Range("D28").Value = "'1.667"
$ws.Range("E28").Value = "  +11.02%  "
$ws.Range("D29").Value = "'7.458"
$ws.Range("E29").Value = "  +2.79%  "
$ws.Range("D30").Value = "'0.05644"
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("E31").Value = "  +3.46%  "
$ws.Range("D32").Value = "'3.706"
$ws.Range("E32").Value = "  +6.15%  "
$ws.Range("D33").Value = "'3.523"
$ws.Range("E33").Value = "  +4.25%  "
$ws.Range("D34").Value = "'1.681"
$ws.Range("E34").Value = "  +6.37%  "
$ws.Range("D35").Value = "'0.9814"
$ws.Range("E35").Value = "  +3.55%  "
$ws.Range("D36").Value = "'2.854"
$ws.Range("E36").Value = "  +2.06%  "
$ws.Range("D37").Value = "'2.446"
$ws.Range("E37").Value = "  +1.50%  "
$ws.Range("D38").Value = "'0.5971"
$ws.Range("E38").Value = "  +3.67%  "
$ws.Range("D39").Value = "'0.01668"
$ws.Range("E39").Value = "  +4.74%  "
$ws.Range("D40").Value = "'5.942"
$ws.Range("E40").Value = "  -0.54%  "
$ws.Range("D41").Value = "'0.8511"
$ws.Range("E41").Value = "  +1.46%  "
$ws.Range("D42").Value = "'1.047.94"
$ws.Range("E42").Value = "  +3.33%  "
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").Value = "'102.07"
$ws.Range("E44").Value = "  +0.43%  "
$ws.Range("D45").Value = "'1.887.99"
$ws.Range("E45").Value = "  +5.41%  "
$ws.Range("D46").Value = "'0.0₈117"
$ws.Range("E46").Value = "  +9.05%  "
$ws.Range("D47").Value = "'59.92"
$ws.Range("D48").Value = "'8.257"
$ws.Range("E48").Value = "  +3.03%  "
$ws.Range("D49").Value = "'0.4437"
$ws.Range("E49").Value = "  +2.15%  "
$ws.Range("D50").Value = "'1.003"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("E51").Value = "  -0.04%  "
